# iphone_price.xlsx - populate the price-list table (rows 2-21) on the active sheet.
# The source feed gives every column as text (including numeric-looking prices and
# review counts such as "113K"); column C is numeric (0) only where the listing has
# no reviews at all. We replicate that by forcing column A/B (and C, by default) to
# Text format before assigning, then switching the genuine-zero rows back to General
# and writing a real number 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Apple iPhone 13", "8249000", "113K", 0),
    @("Apple iPhone 15 Plus", "13499000", "68", 0),
    @("Apple iPhone 15", "11249000", "911", 0),
    @("Apple iPhone 16", "14999000", "0", 1),
    @("Apple iPhone 13", "8249000", "62K", 0),
    @("Apple iPhone 16", "14999000", "42", 0),
    @("Apple iPhone 15", "13749000", "658", 0),
    @("Apple iPhone 13 128 GB", "8249000", "777", 0),
    @("Apple iPhone 16", "14998000", "57", 0),
    @("Apple iPhone 15 128 GB", "11249000", "50", 0),
    @("Apple iPhone 17", "17249000", "236", 0),
    @("Apple iPhone 15", "11249000", "200", 0),
    @("IPhone 17 Pro 17 Pro max 17 Air 16 16E 15 14 plus 13 12 11 Pro Promax X Xs Xr Xsmax New Tempered SUPERFIT Auto Install Anti SPY PRIVACY - HITAM Temperd Super Fit Instal Antigores Gores Screen Protector Tempred Glass Kaca Easy 4g 5g s Liquid Black Gelap", "27256", "0", 1),
    @("Apple iPhone 17 Pro Max", "25749000", "683", 0),
    @("Apple iPhone 16", "14749000", "29", 0),
    @("Apple iPhone 17 Pro", "23749000", "253", 0),
    @("TELEPHONE IP PHONE CISCO 7942G NEW BERGARANSI", "510000", "0", 1),
    @("Apple iPhone 17 Pro", "23749000", "82", 0),
    @("Apple iPhone 14", "9749000", "679", 0),
    @("Apple iPhone 16e", "11749000", "0", 1)
)

$row = 2
foreach ($rec in $data) {
    $name = $rec[0]
    $price = $rec[1]
    $reviews = $rec[2]
    $isNumeric = $rec[3]

    # Column A/B: always text (keeps the leading-digit price string intact).
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $price

    if ($isNumeric -eq 1) {
        # No review count on the listing -> stored as a real numeric 0.
        $ws.Cells.Item($row, 3).NumberFormat = "General"
        $ws.Cells.Item($row, 3).Value = 0
    } else {
        $ws.Cells.Item($row, 3).NumberFormat = "@"
        $ws.Cells.Item($row, 3).Value = $reviews
    }

    $row++
}
